# Updates the benchmark-results table: refreshes the first dozen summary
# rows with new values, and collapses the three trailing multi-column rows
# (each holding a whole tab-separated stats line in a single run) down to
# the single headline figure that now lives in row 1/2/3 respectively.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $text) {
    $c = $table.Cell($row, 1)
    $r = $c.Range
    # Trim the cell's trailing end-of-cell marker before writing new text so
    # we don't leave stray paragraph/cell marks behind.
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $text
}

Set-CellText $t 1  "0M"
Set-CellText $t 2  "0M"
Set-CellText $t 3  "0M"
Set-CellText $t 4  "1909"
Set-CellText $t 5  "0.00001"
Set-CellText $t 6  "0.00272"
Set-CellText $t 7  "0.00018"
Set-CellText $t 8  "0.00008"
Set-CellText $t 9  "0.00032"
Set-CellText $t 10 "0.00038"
Set-CellText $t 11 "0.00050"
Set-CellText $t 12 "0.41839"

Set-CellText $t 44 "99.74"
Set-CellText $t 45 "0.42"
Set-CellText $t 46 "163"
